$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 to make room for the header row
$ws.Rows.Item(1).Insert()

# Fill in the header row with the column names
$ws.Range("A1").Value = "GeoNameId"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Country"
$ws.Range("D1").Value = "Latitude"
$ws.Range("E1").Value = "Longitude"
$ws.Range("F1").Value = "Altitude"

# Latitude/Longitude header cells keep the decimal number style used by
# the data column beneath them
$ws.Range("D1:E1").NumberFormat = "#,###.0000000"

# Update selection to match target (B12 active cell)
$ws.Range("B12").Select()
